# Update "horarios" (schedules) workbook with freshly scraped data.
# Sheet "LP1912"      -> all lines
# Sheet "LP1912-215"  -> only lines whose "Linea" starts with 215
# Sheet "6203-6173"   -> all lines (mirrors "LP1912" data)

$wb = $excel.ActiveWorkbook

$lastUpdate = "06:28:32"

# Full dataset (Hora_Scrap|Hora_Llegada|Linea|Minutos) for the "all lines" sheets.
$allLinesCsv = @"
06:28:32|06:28|23_HERNANDEZ|0
04:01:06|04:02|81_EL PELIGRO|1
04:58:02|04:59|11_ETCHEVERRY|1
04:44:55|04:46|215_EL PELIGRO|2
06:28:32|06:32|11_ETCHEVERRY|4
00:46:06|01:12|215_ALUAR|26
04:44:55|05:11|17_ROMERO|27
04:44:55|05:21|23_HERNANDEZ|37
04:01:06|04:47|215_EL PELIGRO|46
04:44:55|05:31|81_EL PELIGRO|47
06:28:32|07:16|16_SANTA ANA|48
05:55:51|06:50|215A_EL PATO|55
05:55:51|06:53|14_ABASTO|58
04:44:55|05:43|14_ABASTO|59
03:46:12|04:46|215A_EL PATO|60
06:28:32|07:29|14_ABASTO|61
06:28:32|07:29|215A_LA PLATA|61
05:55:51|07:00|16_SANTA ANA|65
01:55:38|03:02|15_ABASTO|67
04:44:55|05:51|17_ROMERO|67
06:28:32|07:36|17X38_ROMERO|68
05:55:51|07:03|225_GOMEZ|68
04:01:06|05:12|17_ROMERO|71
05:55:51|07:06|215C_EL PATO|71
00:46:06|01:58|14_ABASTO|72
04:30:03|05:44|14_ABASTO|74
04:44:55|06:00|16_SANTA ANA|76
05:55:51|07:13|14X44_ABASTO|78
04:44:55|06:03|10_OLMOS|79
05:37:13|07:01|16_SANTA ANA|84
04:44:55|06:10|215A_EL PATO|86
05:37:13|07:04|225_GOMEZ|87
03:46:12|05:16|17_ROMERO|90
05:37:13|07:07|215C_EL PATO|90
06:28:32|07:58|23_HERNANDEZ|90
06:28:32|07:59|11_ETCHEVERRY|91
04:01:06|05:32|81_EL PELIGRO|91
04:30:03|06:01|16_SANTA ANA|91
02:29:13|04:01|81_EL PELIGRO|92
04:58:02|06:31|17X38_ROMERO|93
05:55:51|07:28|14_ABASTO|93
04:58:02|06:31|16_SANTA ANA|93
04:30:03|06:04|10_OLMOS|94
06:28:32|08:03|17_ROMERO|95
01:22:42|02:58|215_ALUAR|96
03:46:12|05:22|23_HERNANDEZ|96
05:37:13|07:14|14X44_ABASTO|97
04:44:55|06:23|11_ETCHEVERRY|99
05:55:51|07:35|17X38_ROMERO|100
06:28:32|08:09|10_OLMOS|101
04:30:03|06:11|215A_EL PATO|101
04:58:02|06:39|225_C ROCA-H SUR|101
05:55:51|07:36|27_EL RETIRO|101
04:44:55|06:26|23_HERNANDEZ|102
06:28:32|08:10|215A_LA PLATA|102
05:37:13|07:21|215A_EL PATO|104
04:01:06|05:45|14_ABASTO|104
04:44:55|06:30|16_SANTA ANA|106
04:44:55|06:30|17X38_ROMERO|106
06:28:32|08:15|17_ROMERO|107
05:55:51|07:43|10_OLMOS|108
03:46:12|05:35|215B_EL PATO|109
05:37:13|07:27|215A_LA PLATA|110
04:01:06|05:52|17_ROMERO|111
04:58:02|06:51|215A_EL PATO|113
03:00:53|04:53|11_ETCHEVERRY|113
01:55:38|03:48|14_ABASTO|113
04:44:55|06:38|225_C ROCA-H SUR|114
04:30:03|06:24|11_ETCHEVERRY|114
06:28:32|08:23|215C_LA PLATA|115
04:58:02|06:54|14_ABASTO|116
05:55:51|07:51|15_ABASTO|116
05:37:13|07:33|23_HERNANDEZ|116
04:30:03|06:27|23_HERNANDEZ|117
02:47:42|04:45|215A_EL PATO|118
06:28:32|08:26|15X38_ABASTO|118
"@

# Subset limited to "215..." lines, used by the "LP1912-215" sheet.
$lines215Csv = @"
04:44:55|04:46|215_EL PELIGRO|2
00:46:06|01:12|215_ALUAR|26
04:01:06|04:47|215_EL PELIGRO|46
05:55:51|06:50|215A_EL PATO|55
03:46:12|04:46|215A_EL PATO|60
06:28:32|07:29|215A_LA PLATA|61
05:55:51|07:06|215C_EL PATO|71
04:44:55|06:10|215A_EL PATO|86
05:37:13|07:07|215C_EL PATO|90
01:22:42|02:58|215_ALUAR|96
04:30:03|06:11|215A_EL PATO|101
06:28:32|08:10|215A_LA PLATA|102
05:37:13|07:21|215A_EL PATO|104
03:46:12|05:35|215B_EL PATO|109
05:37:13|07:27|215A_LA PLATA|110
04:58:02|06:51|215A_EL PATO|113
06:28:32|08:23|215C_LA PLATA|115
02:47:42|04:45|215A_EL PATO|118
"@

function Write-ScheduleSheet {
    param(
        $Worksheet,
        [string]$Csv,
        [string]$LastUpdate
    )

    $dataLines = $Csv -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }

    # Header / metadata block (rows 2-3)
    $Worksheet.Cells.Item(2, 1).Value = "Última actualización: $LastUpdate"
    $Worksheet.Cells.Item(3, 1).Value = "Total filas: $($dataLines.Count)"

    # Data rows start at row 6
    $r = 6
    foreach ($line in $dataLines) {
        $parts = $line -split '\|'
        $Worksheet.Cells.Item($r, 1).Value = $parts[0]
        $Worksheet.Cells.Item($r, 2).Value = $parts[1]
        $Worksheet.Cells.Item($r, 3).Value = $parts[2]
        $Worksheet.Cells.Item($r, 4).Value = [int]$parts[3]
        $r++
    }
}

$wsLP1912    = $wb.Worksheets.Item("LP1912")
$wsLP1912215 = $wb.Worksheets.Item("LP1912-215")
$wsRoca      = $wb.Worksheets.Item("6203-6173")

# NOTE: use positional arguments (named parameter binding misbehaves for
# COM object parameters in this runtime).
Write-ScheduleSheet $wsLP1912    $allLinesCsv $lastUpdate
Write-ScheduleSheet $wsLP1912215 $lines215Csv $lastUpdate
Write-ScheduleSheet $wsRoca      $allLinesCsv $lastUpdate
